$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 458.42105
$ws.Range("I33").Value = 458.42105
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 458.42105
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -229.42105
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1694.1666
$ws.Range("J116").Value = 1830.6
$ws.Range("L116").Value = 1830.6
$ws.Range("N116").Value = -8714.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1000
$ws.Range("I127").Value = 1000
$ws.Range("K127").Value = 3000
$ws.Range("M127").Value = 1960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 19092.564
$ws.Range("I129").Value = 575.6667
$ws.Range("J129").Value = 24260.07
$ws.Range("K129").Value = 1727.0001
$ws.Range("L129").Value = 72780.20999999999
$ws.Range("M129").Value = 3272.9999
$ws.Range("N129").Value = -82780.20999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3073.9177
$ws.Range("I32").Value = 2631.0376
$ws.Range("J32").Value = 10160
$ws.Range("K32").Value = 2631.0376
$ws.Range("L32").Value = 10160
$ws.Range("M32").Value = -2344.0376
$ws.Range("N32").Value = -10734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1157
$ws.Range("I45").Value = 1064.2142
$ws.Range("J45").Value = 1342.5714
$ws.Range("K45").Value = 1064.2142
$ws.Range("L45").Value = 1342.5714
$ws.Range("M45").Value = -687.2141999999999
$ws.Range("N45").Value = -2096.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 668270.4399999999
$ws.Range("I63").Value = 910896.2
$ws.Range("J63").Value = 1049.75
$ws.Range("K63").Value = 910896.2
$ws.Range("L63").Value = 1049.75
$ws.Range("M63").Value = -910210.2
$ws.Range("N63").Value = -2421.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 668270.4399999999
$ws.Range("I66").Value = 910896.2
$ws.Range("J66").Value = 1049.75
$ws.Range("K66").Value = 4554481
$ws.Range("L66").Value = 5248.75
$ws.Range("M66").Value = -4551049
$ws.Range("N66").Value = -12112.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6859.3687
$ws.Range("I110").Value = 8074.3335
$ws.Range("J110").Value = 2303.25
$ws.Range("K110").Value = 8074.3335
$ws.Range("L110").Value = 2303.25
$ws.Range("M110").Value = -6029.3335
$ws.Range("N110").Value = -6393.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2892.1167
$ws.Range("I132").Value = 2985.1956
$ws.Range("J132").Value = 2586.2856
$ws.Range("K132").Value = 8955.586800000001
$ws.Range("L132").Value = 7758.8568
$ws.Range("M132").Value = -6425.586800000001
$ws.Range("N132").Value = -12818.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24134.756
$ws.Range("I134").Value = 36143.17
$ws.Range("J134").Value = 2369.5
$ws.Range("K134").Value = 108429.51
$ws.Range("L134").Value = 7108.5
$ws.Range("M134").Value = -105894.51
$ws.Range("N134").Value = -12178.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2288.7273
$ws.Range("I99").Value = 2058.2222
$ws.Range("J99").Value = 2448.3076
$ws.Range("K99").Value = 2058.2222
$ws.Range("L99").Value = 2448.3076
$ws.Range("M99").Value = -560.2222000000002
$ws.Range("N99").Value = -5444.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 936.6842
$ws.Range("I105").Value = 689.75
$ws.Range("J105").Value = 1360
$ws.Range("K105").Value = 689.75
$ws.Range("L105").Value = 1360
$ws.Range("M105").Value = 1057.25
$ws.Range("N105").Value = -4854

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2288.7273
$ws.Range("I126").Value = 2058.2222
$ws.Range("J126").Value = 2448.3076
$ws.Range("K126").Value = 6174.6666
$ws.Range("L126").Value = 7344.9228
$ws.Range("M126").Value = -3704.6666
$ws.Range("N126").Value = -12284.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1061.1
$ws.Range("J122").Value = 803.75
$ws.Range("L122").Value = 7233.75
$ws.Range("N122").Value = -12133.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2001.8
$ws.Range("I132").Value = 1252
$ws.Range("J132").Value = 2501.6667
$ws.Range("K132").Value = 11268
$ws.Range("L132").Value = 22515.0003
$ws.Range("M132").Value = -8738
$ws.Range("N132").Value = -27575.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18218038
$ws.Range("I70").Value = 22177642
$ws.Range("J70").Value = 3860
$ws.Range("K70").Value = 22177642
$ws.Range("L70").Value = 3860
$ws.Range("M70").Value = -22177372
$ws.Range("N70").Value = -4400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18218038
$ws.Range("I73").Value = 22177642
$ws.Range("J73").Value = 3860
$ws.Range("K73").Value = 22177642
$ws.Range("L73").Value = 3860
$ws.Range("M73").Value = -22176706
$ws.Range("N73").Value = -5732

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2755
$ws.Range("I132").Value = 3222.1428
$ws.Range("J132").Value = 1937.5
$ws.Range("K132").Value = 9666.428400000001
$ws.Range("L132").Value = 5812.5
$ws.Range("M132").Value = -7136.428400000001
$ws.Range("N132").Value = -10872.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3984.4546
$ws.Range("I136").Value = 4517.393
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 13552.179
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -11002.179
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2166.5
$ws.Range("I81").Value = 2249.6667
$ws.Range("J81").Value = 2083.3333
$ws.Range("K81").Value = 4499.3334
$ws.Range("L81").Value = 4166.6666
$ws.Range("M81").Value = -3438.3334
$ws.Range("N81").Value = -6288.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2166.5
$ws.Range("I84").Value = 2249.6667
$ws.Range("J84").Value = 2083.3333
$ws.Range("K84").Value = 22496.667
$ws.Range("L84").Value = 20833.333
$ws.Range("M84").Value = -17192.667
$ws.Range("N84").Value = -31441.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2799.9412
$ws.Range("I132").Value = 2458.9
$ws.Range("J132").Value = 3287.1428
$ws.Range("K132").Value = 7376.700000000001
$ws.Range("L132").Value = 9861.428400000001
$ws.Range("M132").Value = -4846.700000000001
$ws.Range("N132").Value = -14921.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5250.8623
$ws.Range("I136").Value = 5409.8213
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 16229.4639
$ws.Range("L136").Value = 2400
$ws.Range("M136").Value = -13679.4639
$ws.Range("N136").Value = -7500
